$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on D and E columns for rows 2-51 so that numeric-looking
# strings (e.g. "1.00", "7.30") remain stored as text, matching the source data.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '59.924.44'
$ws.Range("E2").Value = '  +2.00%  '
$ws.Range("D3").Value = '3.185.61'
$ws.Range("E3").Value = '  +0.94%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '533.96'
$ws.Range("E5").Value = '  +0.31%  '
$ws.Range("D6").Value = '145.07'
$ws.Range("E6").Value = '  +3.70%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '0.528'
$ws.Range("E8").Value = '  -0.16%  '
$ws.Range("D9").Value = '7.30'
$ws.Range("E9").Value = '  -0.12%  '
$ws.Range("E10").Value = '  +2.15%  '
$ws.Range("D11").Value = '0.429'
$ws.Range("E11").Value = '  +0.01%  '
$ws.Range("D12").Value = '3.736.59'
$ws.Range("E12").Value = '  +1.04%  '
$ws.Range("E13").Value = '  -1.83%  '
$ws.Range("D14").Value = '25.84'
$ws.Range("E14").Value = '  -0.04%  '
$ws.Range("E15").Value = '  +1.40%  '
$ws.Range("D16").Value = '59.992.86'
$ws.Range("E16").Value = '  +2.06%  '
$ws.Range("D17").Value = '3.209.33'
$ws.Range("E17").Value = '  +1.89%  '
$ws.Range("E18").Value = '  +0.37%  '
$ws.Range("D19").Value = '13.20'
$ws.Range("E19").Value = '  +1.43%  '
$ws.Range("D20").Value = '8.19'
$ws.Range("E20").Value = '  -0.27%  '
$ws.Range("D21").Value = '367.08'
$ws.Range("E21").Value = '  -0.94%  '
$ws.Range("E22").Value = '  +0.04%  '
$ws.Range("D23").Value = '0.520'
$ws.Range("E23").Value = '  +0.34%  '
$ws.Range("D24").Value = '69.51'
$ws.Range("E24").Value = '  -0.62%  '
$ws.Range("E25").Value = '  +0.18%  '
$ws.Range("D26").Value = '8.64'
$ws.Range("E26").Value = '  +5.22%  '
$ws.Range("E27").Value = '  +0.10%  '
$ws.Range("D28").Value = '0.0₃0866'
$ws.Range("E28").Value = '  +0.36%  '
$ws.Range("D29").Value = '22.29'
$ws.Range("E29").Value = '  +1.47%  '
$ws.Range("E30").Value = '  +0.44%  '
$ws.Range("D31").Value = '6.04'
$ws.Range("E31").Value = '  -1.23%  '
$ws.Range("D32").Value = '5.27'
$ws.Range("E32").Value = '  +1.23%  '
$ws.Range("D33").Value = '1.19'
$ws.Range("E33").Value = '  +1.75%  '
$ws.Range("D34").Value = '6.56'
$ws.Range("E34").Value = '  +4.42%  '
$ws.Range("D35").Value = '155.53'
$ws.Range("E35").Value = '  -2.30%  '
$ws.Range("D36").Value = '1.35'
$ws.Range("E36").Value = '  +0.08%  '
$ws.Range("D37").Value = '2.812.31'
$ws.Range("E37").Value = '  +6.23%  '
$ws.Range("D38").Value = '26.00'
$ws.Range("E38").Value = '  +2.77%  '
$ws.Range("D39").Value = '0.0704'
$ws.Range("E39").Value = '  +3.30%  '
$ws.Range("D40").Value = '1.66'
$ws.Range("E40").Value = '  -0.86%  '
$ws.Range("D41").Value = '4.21'
$ws.Range("E41").Value = '  -0.38%  '
$ws.Range("D42").Value = '0.0296'
$ws.Range("E42").Value = '  +3.44%  '
$ws.Range("D43").Value = '39.57'
$ws.Range("E43").Value = '  +2.03%  '
$ws.Range("E45").Value = '  +1.60%  '
$ws.Range("D46").Value = '3.230.70'
$ws.Range("E46").Value = '  +1.09%  '
$ws.Range("D47").Value = '0.983'
$ws.Range("E47").Value = '  -0.16%  '
$ws.Range("D48").Value = '6.14'
$ws.Range("E48").Value = '  -0.98%  '
$ws.Range("D49").Value = '0.799'
$ws.Range("E49").Value = '  +5.81%  '
$ws.Range("E50").Value = '  +1.21%  '
$ws.Range("E51").Value = '  +0.12%  '
